# "Generate Report for Handback" -- mark the two localized files (zh-cn,
# de-de) as handed back: update the status text, stamp the handback
# file/datetime columns, and widen the columns that now hold longer text.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$srcMdDisplay = "cb8b6387-ef20-4e9d-963c-e493c86ea27a.md"
$srcMdUrl     = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d20b2bcf0278f146f9651082779c235fb1989893/e2e/cb8b6387-ef20-4e9d-963c-e493c86ea27a.md"
$ffMdDisplay  = "ffff9ba7ca40-387a-4f0e-8e98-85c9db8fb613.md"
$ffMdUrl      = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d20b2bcf0278f146f9651082779c235fb1989893/e2e/ffff9ba7ca40-387a-4f0e-8e98-85c9db8fb613.md"

# ---- Overview sheet: status columns (zh-cn = E, de-de = F) ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527

function Update-LocaleSheet {
    param(
        [string]$SheetName,
        [string]$HandbackFileName,
        [string]$HandbackDateTime
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Status column (C)
    $ws.Range("C2").Value = $statusText
    $ws.Range("C3").Value = $statusText

    # Latest Target File (I), Latest Handback File (J), Latest Handback
    # DateTime (K) for both data rows.
    $ws.Range("I2").Value = $srcMdDisplay
    $ws.Range("J2").Value = $HandbackFileName
    $ws.Range("K2").Value = $HandbackDateTime

    $ws.Range("I3").Value = $srcMdDisplay
    $ws.Range("J3").Value = $HandbackFileName
    $ws.Range("K3").Value = $HandbackDateTime

    # Rebuild the hyperlinks so the new "Latest Target File" links sit next
    # to the existing "Source File Name" links, in row order.
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $srcMdUrl, "", "", $srcMdDisplay)
    $ws.Hyperlinks.Add($ws.Range("I2"), $srcMdUrl, "", "", $srcMdDisplay)
    $ws.Hyperlinks.Add($ws.Range("A3"), $ffMdUrl, "", "", $ffMdDisplay)
    $ws.Hyperlinks.Add($ws.Range("I3"), $srcMdUrl, "", "", $srcMdDisplay)

    # Widen columns that now carry longer values.
    $ws.Columns.Item(3).ColumnWidth = 29.9777047293527
    $ws.Columns.Item(9).ColumnWidth = 40
    $ws.Columns.Item(10).ColumnWidth = 40
}

Update-LocaleSheet "zh-cn" "cb8b6387-ef20-4e9d-963c-e493c86ea27a.379c70dac170965e32cac0c420ccc280e6549c9b.zh-cn.xlf" "2016-08-26 09:08:42"
Update-LocaleSheet "de-de" "cb8b6387-ef20-4e9d-963c-e493c86ea27a.379c70dac170965e32cac0c420ccc280e6549c9b.de-de.xlf" "2016-08-26 09:08:49"
